# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD:AF, styled like the existing header cells ---
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): constant season record for every player row ---
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 70   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 92   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}

Write-Output "Season record columns added."
